$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cnnScriptResults")

# Update values in column B (Test set Accuracy) and C (Time usage)
$ws.Range("B2").Value = 10.1
$ws.Range("B3").Value = 10.9
$ws.Range("B4").Value = 69.8
$ws.Range("B5").Value = 93.2
$ws.Range("C5").Value = 44
$ws.Range("B6").Value = 98.7
$ws.Range("C6").Value = 429

# Update the active selection on the sheet to D6
$ws.Activate()
$ws.Range("D6").Select()
